$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.695.37"
$ws.Range("E2").Value = "'  +2.53%  "
$ws.Range("D3").Value = "'3.190.96"
$ws.Range("E3").Value = "'  +5.73%  "
$ws.Range("D5").Value = "'569.10"
$ws.Range("E5").Value = "'  +3.88%  "
$ws.Range("D6").Value = "'149.38"
$ws.Range("E6").Value = "'  +9.13%  "
$ws.Range("E7").Value = "'  -0.20%  "
$ws.Range("D8").Value = "'3.187.43"
$ws.Range("E8").Value = "'  +5.87%  "
$ws.Range("D9").Value = "'0.507"
$ws.Range("E9").Value = "'  +5.24%  "
$ws.Range("D10").Value = "'6.88"
$ws.Range("E10").Value = "'  +9.07%  "
$ws.Range("D11").Value = "'0.160"
$ws.Range("E11").Value = "'  +6.43%  "
$ws.Range("D12").Value = "'0.481"
$ws.Range("E12").Value = "'  +6.26%  "
$ws.Range("D13").Value = "'38.00"
$ws.Range("E13").Value = "'  +7.95%  "
$ws.Range("D14").Value = "'0.0000230"
$ws.Range("E14").Value = "'  +6.48%  "
$ws.Range("D15").Value = "'3.719.87"
$ws.Range("E15").Value = "'  +6.08%  "
$ws.Range("D16").Value = "'65.889.24"
$ws.Range("E16").Value = "'  +2.87%  "
$ws.Range("D17").Value = "'3.204.86"
$ws.Range("E17").Value = "'  +5.85%  "
$ws.Range("D18").Value = "'532.11"
$ws.Range("E18").Value = "'  +11.87%  "
$ws.Range("E19").Value = "'  +2.96%  "
$ws.Range("D20").Value = "'7.06"
$ws.Range("E20").Value = "'  +8.50%  "
$ws.Range("D21").Value = "'14.43"
$ws.Range("E21").Value = "'  +7.50%  "
$ws.Range("D22").Value = "'0.735"
$ws.Range("E22").Value = "'  +9.11%  "
$ws.Range("D23").Value = "'7.66"
$ws.Range("E23").Value = "'  +9.43%  "
$ws.Range("D24").Value = "'13.34"
$ws.Range("E24").Value = "'  +8.78%  "
$ws.Range("D25").Value = "'80.48"
$ws.Range("E25").Value = "'  +3.96%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "'  +0.13%  "
$ws.Range("D27").Value = "'9.19"
$ws.Range("E27").Value = "'  +20.83%  "
$ws.Range("D28").Value = "'2.90"
$ws.Range("E28").Value = "'  +8.30%  "
$ws.Range("D29").Value = "'2.23"
$ws.Range("E29").Value = "'  +8.97%  "
$ws.Range("D30").Value = "'27.22"
$ws.Range("E30").Value = "'  +7.18%  "
$ws.Range("E31").Value = "'  +0.22%  "
$ws.Range("D32").Value = "'2.69"
$ws.Range("E32").Value = "'  +5.05%  "
$ws.Range("E33").Value = "'  +5.95%  "
$ws.Range("D34").Value = "'551.99"
$ws.Range("E34").Value = "'  -1.23%  "
$ws.Range("D35").Value = "'6.28"
$ws.Range("E35").Value = "'  +9.24%  "
$ws.Range("D36").Value = "'5.56"
$ws.Range("E36").Value = "'  +5.76%  "
$ws.Range("D37").Value = "'54.62"
$ws.Range("E37").Value = "'  +5.44%  "
$ws.Range("D38").Value = "'0.0445"
$ws.Range("E38").Value = "'  +8.78%  "
$ws.Range("D39").Value = "'0.0847"
$ws.Range("E39").Value = "'  +8.27%  "
$ws.Range("D40").Value = "'0.127"
$ws.Range("E40").Value = "'  +7.11%  "
$ws.Range("D41").Value = "'3.179.97"
$ws.Range("E41").Value = "'  +9.71%  "
$ws.Range("D42").Value = "'2.86"
$ws.Range("E42").Value = "'  +5.73%  "
$ws.Range("D43").Value = "'8.52"
$ws.Range("E43").Value = "'  +5.24%  "
$ws.Range("D44").Value = "'0.278"
$ws.Range("E44").Value = "'  +17.00%  "
$ws.Range("D45").Value = "'2.30"
$ws.Range("E45").Value = "'  +12.78%  "
$ws.Range("D46").Value = "'26.47"
$ws.Range("E46").Value = "'  +8.98%  "
$ws.Range("D48").Value = "'0.0₃0542"
$ws.Range("E48").Value = "'  +4.59%  "
$ws.Range("D49").Value = "'123.48"
$ws.Range("E49").Value = "'  +5.16%  "
$ws.Range("D50").Value = "'0.111"
$ws.Range("E50").Value = "'  +4.44%  "
$ws.Range("D51").Value = "'2.18"
$ws.Range("E51").Value = "'  +9.14%  "
